# Updated coremark data after fixing bug in coremk_or0
# (row 16 = "coremk_or0    without AOT or LW optimisations",
#  row 17 = "coremk_or0").  Downstream formulas in F16/F17 and in the
# "RAW DATA" summary block (rows 33-35, which read off rows 16-18)
# recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16 (coremk_or0 without AOT or LW optimisations) ---
$ws.Range("E16").Value = 774.4
$ws.Range("H16").Value = 361.3
$ws.Range("J16").Value = 288.8
$ws.Range("L16").Value = 775

# --- Row 17 (coremk_or0) ---
$ws.Range("E17").Value = 363.6
$ws.Range("H17").Value = 232
$ws.Range("J17").Value = 75.7
$ws.Range("L17").Value = 363.8

# Selection left by the editor after making the change: the summary
# block A31:M43 was selected, with the active cell at its bottom-right
# corner (M43).
$ws.Range("A31:M43").Select()
